$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Organizacija name update
$ws.Range("B4").Value = "Organizacija 1"

# Row 7: Drvo -> Cement
$ws.Range("A7").Value = "Cement"
$ws.Range("B7").Value = 5
$ws.Range("C7").Value = "kg"
$ws.Range("D7").Value = 30.0
$ws.Range("E7").Value = 150.0

# Row 8: Pesak stays, quantities change, unit becomes m3 (new string)
$ws.Range("A8").Value = "Pesak"
$ws.Range("B8").Value = 3
$ws.Range("C8").Value = "m3"
$ws.Range("D8").Value = 15.2
$ws.Range("E8").Value = 45.6

# Row 9: Silikon -> Stiropor
$ws.Range("A9").Value = "Stiropor"
$ws.Range("B9").Value = 1
$ws.Range("C9").Value = "m2"
$ws.Range("D9").Value = 20.1
$ws.Range("E9").Value = 20.1

# Row 10: total amount updated
$ws.Range("D10").Value = "Ukupan iznos"
$ws.Range("E10").Value = 215.7

$wb.Save()
